$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Expe" + bookmark + "rtise" -> single run "Expertise" (bookmark removed)
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Expertise", $true, $false, $false, $false, $false, $true, 1, $false, "Expertise", 2)

# ---------------------------------------------------------------------------
# 2) Split " experience with various programming paradigms and languages."
#    after "pa" and move the _GoBack bookmark there.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("programming pa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng2)
}

# ---------------------------------------------------------------------------
# 3) "Development of Android application. (Complete product)" paragraph:
#    add a numbered tab stop at 720 twips (36pt) and strip direct/character
#    style formatting from the "Complete product" hyperlink run.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Development of Android application", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $para3 = $rng3.Paragraphs(1)
    $para3.TabStops.Add(36, 6)
}

$rng3b = $d.Content
$found3b = $rng3b.Find.Execute("Complete product", $true, $false, $false, $false, $false, $true, 1, $false, "Complete product", 2)

# ---------------------------------------------------------------------------
# 4) "Being responsive to our team managers at Barclays Lithuania." paragraph:
#    add the same numbered tab stop.
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Being responsive to our team managers at Barclays Lithuania", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $para4 = $rng4.Paragraphs(1)
    $para4.TabStops.Add(36, 6)
}

# ---------------------------------------------------------------------------
# 5) "High-Quality content proposal and creation." paragraph:
#    add the same numbered tab stop.
# ---------------------------------------------------------------------------
$rng5 = $d.Content
$found5 = $rng5.Find.Execute("High-Quality content proposal and creation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found5) {
    $para5 = $rng5.Paragraphs(1)
    $para5.TabStops.Add(36, 6)
}
